$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-01 Sunday", "2026-02-02 Monday"),
    @("346÷2=", "964÷8="),
    @("294÷6=", "181÷8="),
    @("396÷3=", "614÷6="),
    @("445÷3=", "962÷8="),
    @("692÷4=", "422÷7="),
    @("279÷4=", "888÷3="),
    @("701÷3=", "195÷3="),
    @("546÷4=", "489÷2="),
    @("362÷3=", "447÷6="),
    @("389÷9=", "138÷3="),
    @("141÷9=", "933÷3="),
    @("611÷6=", "754÷2="),
    @("698÷8=", "284÷3="),
    @("776÷6=", "233÷7="),
    @("904÷2=", "361÷9="),
    @("804÷4=", "800÷5="),
    @("718÷5=", "714÷9="),
    @("706÷5=", "678÷9="),
    @("576÷8=", "955÷2="),
    @("805÷2=", "151÷2="),
    @("877÷3=", "437÷3="),
    @("920÷2=", "828÷5="),
    @("834÷6=", "898÷3="),
    @("406÷5=", "417÷9="),
    @("585÷2=", "759÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
